$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting existing D:K to F:M
$ws.Columns("D:E").Insert(-4161)

# Copy number formats (date / number styles) from column F into new D,E columns
$ws.Range("F5:F102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Range("F5:F102").Copy()
$ws.Range("E5:E102").PasteSpecial(-4122)

# Set the updated cell values for all rows (D:M)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("F7").Value = 43281
$ws.Range("G7").Value = 43190
$ws.Range("H7").Value = 43100
$ws.Range("I7").Value = 43008
$ws.Range("J7").Value = 42916
$ws.Range("K7").Value = 42825
$ws.Range("L7").Value = 42735
$ws.Range("M7").Value = 42643
$ws.Range("D8").Value = 2842100
$ws.Range("E8").Value = 2759900
$ws.Range("F8").Value = 2774900
$ws.Range("G8").Value = 2766100
$ws.Range("H8").Value = 2828800
$ws.Range("I8").Value = 2790900
$ws.Range("J8").Value = 2742500
$ws.Range("K8").Value = 2627200
$ws.Range("L8").Value = 2697500
$ws.Range("M8").Value = 2652800
$ws.Range("D9").Value = 1435400
$ws.Range("E9").Value = 1308100
$ws.Range("F9").Value = 1359800
$ws.Range("G9").Value = 1287300
$ws.Range("H9").Value = 1307500
$ws.Range("I9").Value = 1306400
$ws.Range("J9").Value = 1286600
$ws.Range("K9").Value = 1254900
$ws.Range("L9").Value = 1213900
$ws.Range("M9").Value = 1308100
$ws.Range("D10").Value = 1406700
$ws.Range("E10").Value = 1451800
$ws.Range("F10").Value = 1415100
$ws.Range("G10").Value = 1478800
$ws.Range("H10").Value = 1521200
$ws.Range("I10").Value = 1484500
$ws.Range("J10").Value = 1455900
$ws.Range("K10").Value = 1372300
$ws.Range("L10").Value = 1483600
$ws.Range("M10").Value = 1344700
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("L12").Value = "NA"
$ws.Range("M12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("D14").Value = -54800
$ws.Range("E14").Value = -388000
$ws.Range("F14").Value = -550200
$ws.Range("G14").Value = -22700
$ws.Range("H14").Value = -21800
$ws.Range("I14").Value = "NA"
$ws.Range("J14").Value = "NA"
$ws.Range("K14").Value = "NA"
$ws.Range("L14").Value = "NA"
$ws.Range("M14").Value = "NA"
$ws.Range("D15").Value = 124800
$ws.Range("E15").Value = 121700
$ws.Range("F15").Value = 118900
$ws.Range("G15").Value = 116700
$ws.Range("H15").Value = 114300
$ws.Range("I15").Value = 121600
$ws.Range("J15").Value = 121200
$ws.Range("K15").Value = 117400
$ws.Range("L15").Value = 115000
$ws.Range("M15").Value = 119200
$ws.Range("D17").Value = 2400500
$ws.Range("E17").Value = 2057300
$ws.Range("F17").Value = 1956900
$ws.Range("G17").Value = 2305400
$ws.Range("H17").Value = 2373400
$ws.Range("I17").Value = 2351300
$ws.Range("J17").Value = 2339600
$ws.Range("K17").Value = 2237400
$ws.Range("L17").Value = 2248900
$ws.Range("M17").Value = 2296100
$ws.Range("D18").Value = 441600
$ws.Range("E18").Value = 702600
$ws.Range("F18").Value = 818000
$ws.Range("G18").Value = 460700
$ws.Range("H18").Value = 455300
$ws.Range("I18").Value = 439700
$ws.Range("J18").Value = 403000
$ws.Range("K18").Value = 389800
$ws.Range("L18").Value = 448600
$ws.Range("M18").Value = 356700
$ws.Range("D20").Value = 1300
$ws.Range("E20").Value = 197200
$ws.Range("F20").Value = 404800
$ws.Range("G20").Value = -4000
$ws.Range("H20").Value = -5300
$ws.Range("I20").Value = 18000
$ws.Range("J20").Value = -6100
$ws.Range("K20").Value = 900
$ws.Range("L20").Value = -23600
$ws.Range("M20").Value = -900
$ws.Range("D21").Value = 1038400
$ws.Range("E21").Value = 1421400
$ws.Range("F21").Value = 1738900
$ws.Range("G21").Value = 969000
$ws.Range("H21").Value = 960500
$ws.Range("I21").Value = 960700
$ws.Range("J21").Value = 898700
$ws.Range("K21").Value = 872800
$ws.Range("L21").Value = 875400
$ws.Range("M21").Value = 845200
$ws.Range("D22").Value = 33300
$ws.Range("E22").Value = 29600
$ws.Range("F22").Value = 27700
$ws.Range("G22").Value = 40300
$ws.Range("H22").Value = 40300
$ws.Range("I22").Value = 61700
$ws.Range("J22").Value = 61700
$ws.Range("K22").Value = 72900
$ws.Range("L22").Value = 54700
$ws.Range("M22").Value = 72700
$ws.Range("D23").Value = 409600
$ws.Range("E23").Value = 870200
$ws.Range("F23").Value = 1195200
$ws.Range("G23").Value = 416400
$ws.Range("H23").Value = 409700
$ws.Range("I23").Value = 396000
$ws.Range("J23").Value = 335200
$ws.Range("K23").Value = 317700
$ws.Range("L23").Value = 370400
$ws.Range("M23").Value = 283200
$ws.Range("D24").Value = 28500
$ws.Range("E24").Value = 55600
$ws.Range("F24").Value = 383400
$ws.Range("G24").Value = 134900
$ws.Range("H24").Value = 20800
$ws.Range("I24").Value = 82500
$ws.Range("J24").Value = 111400
$ws.Range("K24").Value = 70600
$ws.Range("L24").Value = 69000
$ws.Range("M24").Value = 46900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("D26").Value = 381200
$ws.Range("E26").Value = 814600
$ws.Range("F26").Value = 811800
$ws.Range("G26").Value = 281500
$ws.Range("H26").Value = 388900
$ws.Range("I26").Value = 313500
$ws.Range("J26").Value = 223800
$ws.Range("K26").Value = 247100
$ws.Range("L26").Value = 301400
$ws.Range("M26").Value = 236300
$ws.Range("D27").Value = 381200
$ws.Range("E27").Value = 814600
$ws.Range("F27").Value = 811800
$ws.Range("G27").Value = 281500
$ws.Range("H27").Value = 388900
$ws.Range("I27").Value = 313500
$ws.Range("J27").Value = 223800
$ws.Range("K27").Value = 247100
$ws.Range("L27").Value = 301400
$ws.Range("M27").Value = 236300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("D32").Value = -1300
$ws.Range("E32").Value = -197200
$ws.Range("F32").Value = -404800
$ws.Range("G32").Value = 4000
$ws.Range("H32").Value = 5300
$ws.Range("I32").Value = -18000
$ws.Range("J32").Value = 6100
$ws.Range("K32").Value = -900
$ws.Range("L32").Value = 23600
$ws.Range("M32").Value = 900
$ws.Range("D33").Value = 381200
$ws.Range("E33").Value = 814600
$ws.Range("F33").Value = 811800
$ws.Range("G33").Value = 281500
$ws.Range("H33").Value = 388900
$ws.Range("I33").Value = 313500
$ws.Range("J33").Value = 223800
$ws.Range("K33").Value = 247100
$ws.Range("L33").Value = 301400
$ws.Range("M33").Value = 236300
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("D35").Value = 381200
$ws.Range("E35").Value = 814600
$ws.Range("F35").Value = 811800
$ws.Range("G35").Value = 281500
$ws.Range("H35").Value = 388900
$ws.Range("I35").Value = 313500
$ws.Range("J35").Value = 223800
$ws.Range("K35").Value = 247100
$ws.Range("L35").Value = 301400
$ws.Range("M35").Value = 236300
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("F38").Value = 43281
$ws.Range("G38").Value = 43190
$ws.Range("H38").Value = 43100
$ws.Range("I38").Value = 43008
$ws.Range("J38").Value = 42916
$ws.Range("K38").Value = 42825
$ws.Range("L38").Value = 42735
$ws.Range("M38").Value = 42643
$ws.Range("D41").Value = 866900
$ws.Range("E41").Value = 952100
$ws.Range("F41").Value = 1135700
$ws.Range("G41").Value = 1116400
$ws.Range("H41").Value = 1038400
$ws.Range("I41").Value = 1428300
$ws.Range("J41").Value = 1909300
$ws.Range("K41").Value = 1559200
$ws.Range("L41").Value = 1266500
$ws.Range("M41").Value = 1437500
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("D43").Value = 3427300
$ws.Range("E43").Value = 3163800
$ws.Range("F43").Value = 2849800
$ws.Range("G43").Value = 2958400
$ws.Range("H43").Value = 2901400
$ws.Range("I43").Value = 2974400
$ws.Range("J43").Value = 3036200
$ws.Range("K43").Value = 2869600
$ws.Range("L43").Value = 2971200
$ws.Range("M43").Value = 2753100
$ws.Range("D44").Value = 118500
$ws.Range("E44").Value = 117800
$ws.Range("F44").Value = 126000
$ws.Range("G44").Value = 130900
$ws.Range("H44").Value = 89400
$ws.Range("I44").Value = 97600
$ws.Range("J44").Value = 100000
$ws.Range("K44").Value = 98600
$ws.Range("L44").Value = 101800
$ws.Range("M44").Value = 124100
$ws.Range("D45").Value = 295300
$ws.Range("E45").Value = 342700
$ws.Range("F45").Value = 416100
$ws.Range("G45").Value = 386100
$ws.Range("H45").Value = 260400
$ws.Range("I45").Value = 286100
$ws.Range("J45").Value = 402900
$ws.Range("K45").Value = 438400
$ws.Range("L45").Value = 224900
$ws.Range("M45").Value = 258600
$ws.Range("D46").Value = 4707900
$ws.Range("E46").Value = 4576400
$ws.Range("F46").Value = 4527600
$ws.Range("G46").Value = 4591900
$ws.Range("H46").Value = 4289700
$ws.Range("I46").Value = 4786400
$ws.Range("J46").Value = 5448300
$ws.Range("K46").Value = 4965900
$ws.Range("L46").Value = 4564400
$ws.Range("M46").Value = 4573200
$ws.Range("D47").Value = 963000
$ws.Range("E47").Value = 1430500
$ws.Range("F47").Value = 1305900
$ws.Range("G47").Value = 301900
$ws.Range("H47").Value = 286800
$ws.Range("I47").Value = 209000
$ws.Range("J47").Value = 211900
$ws.Range("K47").Value = 215200
$ws.Range("L47").Value = 216000
$ws.Range("M47").Value = 269200
$ws.Range("D48").Value = 8746500
$ws.Range("E48").Value = 8751600
$ws.Range("F48").Value = 8581800
$ws.Range("G48").Value = 8489800
$ws.Range("H48").Value = 8517500
$ws.Range("I48").Value = 8311100
$ws.Range("J48").Value = 8155700
$ws.Range("K48").Value = 7857500
$ws.Range("L48").Value = 7919900
$ws.Range("M48").Value = 7621600
$ws.Range("D49").Value = 10824600
$ws.Range("E49").Value = 10874200
$ws.Range("F49").Value = 10953400
$ws.Range("G49").Value = 11020200
$ws.Range("H49").Value = 11109400
$ws.Range("I49").Value = 11126200
$ws.Range("J49").Value = 11225000
$ws.Range("K49").Value = 10935700
$ws.Range("L49").Value = 11035500
$ws.Range("M49").Value = 11102400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("D52").Value = 1052600
$ws.Range("E52").Value = 1116500
$ws.Range("F52").Value = 1459200
$ws.Range("G52").Value = 1819100
$ws.Range("H52").Value = 1789100
$ws.Range("I52").Value = 1791700
$ws.Range("J52").Value = 1714000
$ws.Range("K52").Value = 1666100
$ws.Range("L52").Value = 1584700
$ws.Range("M52").Value = 1665400
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("D54").Value = 26294700
$ws.Range("E54").Value = 26749200
$ws.Range("F54").Value = 26827900
$ws.Range("G54").Value = 26222800
$ws.Range("H54").Value = 25992500
$ws.Range("I54").Value = 26224400
$ws.Range("J54").Value = 26754900
$ws.Range("K54").Value = 25640400
$ws.Range("L54").Value = 25320600
$ws.Range("M54").Value = 25231800
$ws.Range("D57").Value = 1959500
$ws.Range("E57").Value = 2079200
$ws.Range("F57").Value = 2056500
$ws.Range("G57").Value = 1868600
$ws.Range("H57").Value = 1909300
$ws.Range("I57").Value = 1869200
$ws.Range("J57").Value = 1822100
$ws.Range("K57").Value = 1699100
$ws.Range("L57").Value = 1888200
$ws.Range("M57").Value = 1833200
$ws.Range("D58").Value = 375400
$ws.Range("E58").Value = 387200
$ws.Range("F58").Value = 394100
$ws.Range("G58").Value = 722700
$ws.Range("H58").Value = 777700
$ws.Range("I58").Value = 873800
$ws.Range("J58").Value = 1420400
$ws.Range("K58").Value = 1105200
$ws.Range("L58").Value = 1156900
$ws.Range("M58").Value = 921800
$ws.Range("D59").Value = 2064900
$ws.Range("E59").Value = 2587700
$ws.Range("F59").Value = 2561700
$ws.Range("G59").Value = 1950100
$ws.Range("H59").Value = 1892600
$ws.Range("I59").Value = 2249700
$ws.Range("J59").Value = 2450000
$ws.Range("K59").Value = 1875100
$ws.Range("L59").Value = 2025300
$ws.Range("M59").Value = 2310100
$ws.Range("D60").Value = 4399700
$ws.Range("E60").Value = 5054100
$ws.Range("F60").Value = 5012300
$ws.Range("G60").Value = 4541400
$ws.Range("H60").Value = 4579600
$ws.Range("I60").Value = 4992600
$ws.Range("J60").Value = 5692600
$ws.Range("K60").Value = 4679400
$ws.Range("L60").Value = 5070400
$ws.Range("M60").Value = 5065100
$ws.Range("D61").Value = 1198600
$ws.Range("E61").Value = 1229500
$ws.Range("F61").Value = 1273900
$ws.Range("G61").Value = 1335800
$ws.Range("H61").Value = 1391700
$ws.Range("I61").Value = 1191800
$ws.Range("J61").Value = 1251600
$ws.Range("K61").Value = 1579700
$ws.Range("L61").Value = 1131400
$ws.Range("M61").Value = 1272300
$ws.Range("D62").Value = 2337700
$ws.Range("E62").Value = 2144400
$ws.Range("F62").Value = 2317800
$ws.Range("G62").Value = 2290600
$ws.Range("H62").Value = 2212700
$ws.Range("I62").Value = 2218400
$ws.Range("J62").Value = 2208400
$ws.Range("K62").Value = 2069800
$ws.Range("L62").Value = 1940600
$ws.Range("M62").Value = 1860400
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("D66").Value = 7936100
$ws.Range("E66").Value = 8427900
$ws.Range("F66").Value = 8603900
$ws.Range("G66").Value = 8167800
$ws.Range("H66").Value = 8184000
$ws.Range("I66").Value = 8402800
$ws.Range("J66").Value = 9152500
$ws.Range("K66").Value = 8329000
$ws.Range("L66").Value = 8142400
$ws.Range("M66").Value = 8197700
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
$ws.Range("D72").Value = 1442300
$ws.Range("E72").Value = 2036700
$ws.Range("F72").Value = 1940000
$ws.Range("G72").Value = 1211100
$ws.Range("H72").Value = 965200
$ws.Range("I72").Value = 1541100
$ws.Range("J72").Value = 1320900
$ws.Range("K72").Value = 1083900
$ws.Range("L72").Value = 951500
$ws.Range("M72").Value = 1282900
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("D76").Value = 18358600
$ws.Range("E76").Value = 18321300
$ws.Range("F76").Value = 18223900
$ws.Range("G76").Value = 18055000
$ws.Range("H76").Value = 17808500
$ws.Range("I76").Value = 17821500
$ws.Range("J76").Value = 17602400
$ws.Range("K76").Value = 17311400
$ws.Range("L76").Value = 17178200
$ws.Range("M76").Value = 17034100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("F80").Value = 43281
$ws.Range("G80").Value = 43190
$ws.Range("H80").Value = 43100
$ws.Range("I80").Value = 43008
$ws.Range("J80").Value = 42916
$ws.Range("K80").Value = 42825
$ws.Range("L80").Value = 42735
$ws.Range("M80").Value = 42643
$ws.Range("D81").Value = 381200
$ws.Range("E81").Value = 814600
$ws.Range("F81").Value = 811800
$ws.Range("G81").Value = 281500
$ws.Range("H81").Value = 388900
$ws.Range("I81").Value = 313500
$ws.Range("J81").Value = 223800
$ws.Range("K81").Value = 247100
$ws.Range("L81").Value = 301400
$ws.Range("M81").Value = 236300
$ws.Range("D83").Value = 595500
$ws.Range("E83").Value = 521700
$ws.Range("F83").Value = 516100
$ws.Range("G83").Value = 512300
$ws.Range("H83").Value = 510400
$ws.Range("I83").Value = 503000
$ws.Range("J83").Value = 501800
$ws.Range("K83").Value = 482200
$ws.Range("L83").Value = 450300
$ws.Range("M83").Value = 489300
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
$ws.Range("D89").Value = 1177000
$ws.Range("E89").Value = 625500
$ws.Range("F89").Value = 503800
$ws.Range("G89").Value = 755300
$ws.Range("H89").Value = 759700
$ws.Range("I89").Value = 995800
$ws.Range("J89").Value = 1056700
$ws.Range("K89").Value = 414900
$ws.Range("L89").Value = 823100
$ws.Range("M89").Value = 764700
$ws.Range("D91").Value = -672600
$ws.Range("E91").Value = -555600
$ws.Range("F91").Value = -417300
$ws.Range("G91").Value = -538200
$ws.Range("H91").Value = -721700
$ws.Range("I91").Value = -525300
$ws.Range("J91").Value = -437900
$ws.Range("K91").Value = -445500
$ws.Range("L91").Value = -482400
$ws.Range("M91").Value = -465500
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("D94").Value = -638700
$ws.Range("E94").Value = -208600
$ws.Range("F94").Value = -68800
$ws.Range("G94").Value = -539200
$ws.Range("H94").Value = -705700
$ws.Range("I94").Value = -552100
$ws.Range("J94").Value = -411400
$ws.Range("K94").Value = -478200
$ws.Range("L94").Value = -489300
$ws.Range("M94").Value = -476000
$ws.Range("D96").Value = -547000
$ws.Range("E96").Value = -513400
$ws.Range("F96").Value = -100
$ws.Range("G96").Value = -100
$ws.Range("H96").Value = -606200
$ws.Range("I96").Value = -334200
$ws.Range("J96").Value = -100
$ws.Range("K96").Value = -100
$ws.Range("L96").Value = -470300
$ws.Range("M96").Value = -265400
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("D100").Value = -623500
$ws.Range("E100").Value = -600400
$ws.Range("F100").Value = -415700
$ws.Range("G100").Value = -138100
$ws.Range("H100").Value = -611900
$ws.Range("I100").Value = -924700
$ws.Range("J100").Value = -179400
$ws.Range("K100").Value = 356100
$ws.Range("L100").Value = -504800
$ws.Range("M100").Value = -259100
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 0
$ws.Range("D102").Value = -85200
$ws.Range("E102").Value = -183600
$ws.Range("F102").Value = 19300
$ws.Range("G102").Value = 78000
$ws.Range("H102").Value = -389900
$ws.Range("I102").Value = -481000
$ws.Range("J102").Value = 297900
$ws.Range("K102").Value = 292700
$ws.Range("L102").Value = -171000
$ws.Range("M102").Value = 29500
